$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "249.63"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.89"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.425"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05637"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.419"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.364"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8114"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9182"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1439"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07527"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03106"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03110"
$ws.Range("D13").Style = "Normal"

$ws.Range("B14").Value = "BitMartToken"

$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09340"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"

$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.553"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"

$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001595"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"

$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04753"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"

$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005790"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006391"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005002"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001031"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001499"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.698"
$ws.Range("D23").Style = "Normal"

$ws.Range("B25").Value = "BitpandaEcosystemToken"

$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3299"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"

$ws.Range("B26").Value = "ProBitToken"

$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1296"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "25ProBitTokenPROB"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04037"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006807"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1070"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002719"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007501"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005466"
$ws.Range("D45").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5001"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2408"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
